# "Co sua giao dien" - fix/clean up the sheet interface:
#  - correct the training-program name/code on row 5
#  - drop the leftover PCCM/LT/TH helper columns (K:L) that were never
#    meant to be shown, including the per-row computed hour splits and
#    the K62 grand-total formula
#  - reset the view back to a sane scroll position / selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Update the training program / major line.
$ws.Range("A5").Value = "Ngành đào tạo: Hệ thống thông tin    Mã ngành: 6320201"

# 2) Remove the stray K/L helper columns content.
# K9:L9 keeps its (merged) cell/style but loses the "PCCM" label.
$ws.Range("K9:L9").ClearContents()

# K10:L10 ("LT" / "TH" headers) are removed entirely.
$ws.Range("K10:L10").Clear()

# L11 held a helper SUM formula - remove entirely.
$ws.Range("L11").Clear()

# K16:L18, K29:L31 held per-row LT/TH hour splits - remove entirely.
$ws.Range("K16:L18").Clear()
$ws.Range("K29:L31").Clear()

# K62 held the grand-total helper formula - remove entirely.
$ws.Range("K62").Clear()

# 3) Reset the view: selection back on the title row, no stale scroll offset.
$ws.Range("A6:J6").Select()
